$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# The "Recorded By" column (G) sometimes lists both the user and "System"
# as recorders for a session. Swap the order so "System" is listed first,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

$col = $ws.Range("G1:G319")
[void]$col.Replace($oldVal, $newVal)
